$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.300.17'
$ws.Range("E2").Value = '  +2.05%  '
$ws.Range("D3").Value = '2.098.92'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("E4").Value = '  -0.76%  '
$ws.Range("D5").Value = '''342.69'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  -0.70%  '
$ws.Range("D7").Value = '''0.5284'
$ws.Range("E7").Value = '  +2.24%  '
$ws.Range("D8").Value = '''0.4395'
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").Value = '''55.07'
$ws.Range("E9").Value = '  +2.60%  '
$ws.Range("D10").Value = '''0.09364'
$ws.Range("E10").Value = '  +2.05%  '
$ws.Range("D11").Value = '''1.176'
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").Value = '''24.80'
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("E13").Value = '  +5.23%  '
$ws.Range("D14").Value = '''6.869'
$ws.Range("E14").Value = '  +1.63%  '
$ws.Range("D15").Value = '2.003.51'
$ws.Range("E15").Value = '  -3.66%  '
$ws.Range("D16").Value = '''101.20'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").Value = '''0.00001158'
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").Value = '''21.15'
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("D20").Value = '''0.06727'
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("D21").Value = '''6.383'
$ws.Range("E21").Value = '  +3.03%  '
$ws.Range("D22").Value = '''1.001'
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("D23").Value = '30.286.48'
$ws.Range("E23").Value = '  +1.82%  '
$ws.Range("D24").Value = '''12.43'
$ws.Range("E24").Value = '  -1.73%  '
$ws.Range("D25").Value = '''2.319'
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("D26").Value = '''6.997'
$ws.Range("E26").Value = '  +11.03%  '
$ws.Range("D27").Value = '''21.84'
$ws.Range("E27").Value = '  -0.19%  '
$ws.Range("D28").Value = '''162.68'
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("E29").Value = '  +1.18%  '
$ws.Range("D30").Value = '''133.84'
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("D31").Value = '''1.136'
$ws.Range("E31").Value = '  +0.72%  '
$ws.Range("D32").Value = '''1.690'
$ws.Range("E32").Value = '  +1.41%  '
$ws.Range("E33").Value = '  +0.37%  '
$ws.Range("D34").Value = '''6.258'
$ws.Range("E34").Value = '  +1.21%  '
$ws.Range("D35").Value = '''3.876'
$ws.Range("E35").Value = '  -2.18%  '
$ws.Range("D36").Value = '''10.07'
$ws.Range("E36").Value = '  -3.34%  '
$ws.Range("D37").Value = '''0.02623'
$ws.Range("E37").Value = '  +1.84%  '
$ws.Range("D38").Value = '''0.06766'
$ws.Range("E38").Value = '  +1.11%  '
$ws.Range("E39").Value = '  +1.98%  '
$ws.Range("D40").Value = '''1.353'
$ws.Range("E40").Value = '  +1.53%  '
$ws.Range("D41").Value = '''0.6970'
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("D42").Value = '''0.2218'
$ws.Range("E42").Value = '  +0.34%  '
$ws.Range("D43").Value = '''0.6783'
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").Value = '''14.35'
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").Value = '''2.335'
$ws.Range("E45").Value = '  +0.79%  '
$ws.Range("E46").Value = '  -0.64%  '
$ws.Range("D47").Value = '''1.310'
$ws.Range("E47").Value = '  +8.52%  '
$ws.Range("D48").Value = '''3.642'
$ws.Range("E48").Value = '  +0.56%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '''0.00000000349'
$ws.Range("E49").Value = '  -3.36%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '''1.212'
$ws.Range("E50").Value = '  +6.48%  '
$ws.Range("D51").Value = '''0.07302'
$ws.Range("E51").Value = '  +3.53%  '
